# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The "Periodo Mora" / "Valor Mora" figures for the two worker rows (16 and
# 17) in the EC table were corrected: the period values 1811/1810 and their
# matching mora amounts 40000/20000 were swapped between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 ("1811" / 40000)  -> becomes "1810" / 20000
# Row 17 ("1810" / 20000)  -> becomes "1811" / 40000
$ws.Range("E16").Value = "1810"
$ws.Range("F16").Value = 20000

$ws.Range("E17").Value = "1811"
$ws.Range("F17").Value = 40000
